$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '48.278.99'
$ws.Range("E2").Value = '  +2.07%  '
$ws.Range("D3").Value = '2.526.52'
$ws.Range("E3").Value = '  +1.10%  '
$ws.Range("E4").Value = '  +0.04%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '323.68'
$c.ClearFormats()
$ws.Range("E5").Value = '  +0.16%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '109.56'
$c.ClearFormats()
$ws.Range("E6").Value = '  +0.47%  '
$ws.Range("E7").Value = '  +0.86%  '
$ws.Range("E8").Value = '  +0.01%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.559'
$c.ClearFormats()
$ws.Range("E9").Value = '  +4.35%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '40.76'
$c.ClearFormats()
$ws.Range("E10").Value = '  +4.12%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '20.44'
$c.ClearFormats()
$ws.Range("E11").Value = '  +11.46%  '
$ws.Range("E12").Value = '  +1.40%  '
$ws.Range("E13").Value = '  +0.99%  '
$ws.Range("E14").Value = '  +1.60%  '
$ws.Range("D15").Value = '2.923.75'
$ws.Range("E15").Value = '  +1.13%  '
$ws.Range("D16").Value = '2.529.22'
$ws.Range("E16").Value = '  +0.93%  '
$ws.Range("E17").Value = '  +0.99%  '
$ws.Range("D18").Value = '48.131.03'
$ws.Range("E18").Value = '  +1.91%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '13.31'
$c.ClearFormats()
$ws.Range("E19").Value = '  +3.99%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '6.65'
$c.ClearFormats()
$ws.Range("E20").Value = '  +0.17%  '
$ws.Range("D21").Value = '0.0₃0951'
$ws.Range("E21").Value = '  +1.05%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '2.71'
$c.ClearFormats()
$ws.Range("E22").Value = '  +0.47%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '72.57'
$c.ClearFormats()
$ws.Range("E23").Value = '  +2.78%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '271.19'
$c.ClearFormats()
$ws.Range("E24").Value = '  +9.42%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '2.59'
$c.ClearFormats()
$ws.Range("E25").Value = '  -0.56%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '26.28'
$c.ClearFormats()
$ws.Range("E26").Value = '  +0.83%  '
$ws.Range("E27").Value = '  +0.08%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '10.18'
$c.ClearFormats()
$ws.Range("E28").Value = '  +1.24%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '0.147'
$c.ClearFormats()
$ws.Range("E29").Value = '  +6.20%  '
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '2.21'
$c.ClearFormats()
$ws.Range("E30").Value = '  -3.91%  '
$ws.Range("B31").Value = 'InjectiveProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '35.93'
$c.ClearFormats()
$ws.Range("E31").Value = '  +1.84%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '49.83'
$c.ClearFormats()
$ws.Range("E32").Value = '  -0.01%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '20.01'
$c.ClearFormats()
$ws.Range("E33").Value = '  +0.07%  '
$ws.Range("E34").Value = '  -0.18%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '0.0792'
$c.ClearFormats()
$ws.Range("E36").Value = '  +0.44%  '
$ws.Range("E37").Value = '  +1.04%  '
$ws.Range("E38").Value = '  +1.44%  '
$ws.Range("E39").Value = '  +0.30%  '
$ws.Range("E40").Value = '  +0.00%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '22.34'
$c.ClearFormats()
$ws.Range("E41").Value = '  +4.98%  '
$ws.Range("E42").Value = '  -2.07%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '118.87'
$c.ClearFormats()
$ws.Range("E43").Value = '  -2.05%  '
$ws.Range("E44").Value = '  +0.52%  '
$ws.Range("D45").Value = '2.016.33'
$ws.Range("E45").Value = '  +1.31%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '3.15'
$c.ClearFormats()
$ws.Range("E46").Value = '  +2.91%  '
$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '2.05'
$c.ClearFormats()
$ws.Range("E47").Value = '  -0.48%  '
$ws.Range("B48").Value = 'Stacks'
$ws.Range("C48").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '1.88'
$c.ClearFormats()
$ws.Range("E48").Value = '  +5.64%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '9.16'
$c.ClearFormats()
$ws.Range("E49").Value = '  +0.78%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '5.27'
$c.ClearFormats()
$ws.Range("E50").Value = '  +1.01%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '80.39'
$c.ClearFormats()
$ws.Range("E51").Value = '  +3.10%  '
